# Final commit of upload excel file
# Updates the "contacts" sheet with corrected/updated values and
# refreshes the related number-formatted cells + row heights.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Rohan / Tintu) ---
$ws.Range("B2").Value = "rohan"
$ws.Range("J2").Value = "tintu@gmail.com"
$ws.Range("L2").Value = "Reading ,Drawing"

# --- Row 3 (Mini) ---
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# --- Font color refresh for the Pincode / Phone columns ---
$ws.Range("I2:I3").Font.Color = 0
$ws.Range("K2:K3").Font.Color = 0

# --- Row heights bumped slightly after the re-entry ---
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5
